$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 209, shifting existing rows 209-216 down to 210-217
$ws.Rows.Item(209).Insert()

# Populate the newly inserted row 209 with values (mirrors the row that used
# to be at 209, except for the updated measurement columns)
$ws.Cells.Item(209, 1).Value = 3
$ws.Cells.Item(209, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(209, 3).Value = "Coquimbo"
$ws.Cells.Item(209, 4).Value = 44509
$ws.Cells.Item(209, 5).Value = 5
$ws.Cells.Item(209, 6).Value = 100112012
$ws.Cells.Item(209, 7).Value = "Espinaca"
$ws.Cells.Item(209, 8).Value = "Sin especificar"
$ws.Cells.Item(209, 9).Value = "Primera"
$ws.Cells.Item(209, 10).Value = 290
$ws.Cells.Item(209, 11).Value = 2000
$ws.Cells.Item(209, 12).Value = 2300
$ws.Cells.Item(209, 13).Value = 2166
$ws.Cells.Item(209, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(209, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(209, 16).Value = 722
$ws.Cells.Item(209, 17).Value = 3
$ws.Cells.Item(209, 18).Value = "Hortaliza"

# Match the date-time number format used by the rest of column D
$ws.Cells.Item(209, 4).NumberFormat = $ws.Cells.Item(210, 4).NumberFormat
